$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts all existing columns
# (A:V) right by one (to B:W) and keeps their formatting, merged cells, etc.
$ws.Columns.Item(1).Insert()

# Header label for the newly inserted "Match ID" column
$ws.Range("A2").Value = "Match ID"

# Bold the header + visible/hidden data cells of the new column
# (this creates the new cell style used by the Match ID column)
$ws.Range("A2:A20").Font.Bold = $true

# Fill the new "Match ID" column with the constant match id value
$ws.Range("A4:A20").Value = 24
$ws.Cells.Item(21, 1).Value = 24

# Row 21 is hidden; re-autofit it so no explicit custom row height sticks
# after writing the value (keeps the row definition clean/unchanged)
$ws.Rows.Item(21).AutoFit()

# Update the selection to reflect the new "Match ID" column selection
$ws.Range("A2:A20").Select()
